$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H31").Value = 100
$ws.Range("I31").Value = 100
$ws.Range("J31").Value = 0
$ws.Range("K31").Value = 300
$ws.Range("L31").Value = 0
$ws.Range("M31").Value = -70
$ws.Range("H41").Value = 737.5
$ws.Range("I41").Value = 316.66666
$ws.Range("J41").Value = 2000
$ws.Range("K41").Value = 316.66666
$ws.Range("L41").Value = 2000
$ws.Range("M41").Value = 123.33334
$ws.Range("N41").Value = -2880
$ws.Range("H53").Value = 1656.25
$ws.Range("I53").Value = 3363.3333
$ws.Range("J53").Value = 632
$ws.Range("K53").Value = 3363.3333
$ws.Range("L53").Value = 632
$ws.Range("M53").Value = -2726.3333
$ws.Range("N53").Value = -1906
$ws.Range("H70").Value = 2570.7727
$ws.Range("I70").Value = 2795.2
$ws.Range("J70").Value = 2383.75
$ws.Range("K70").Value = 8385.599999999999
$ws.Range("L70").Value = 7151.25
$ws.Range("M70").Value = -8115.599999999999
$ws.Range("N70").Value = -7691.25
$ws.Range("H73").Value = 2570.7727
$ws.Range("I73").Value = 2795.2
$ws.Range("J73").Value = 2383.75
$ws.Range("K73").Value = 8385.599999999999
$ws.Range("L73").Value = 7151.25
$ws.Range("M73").Value = -7449.599999999999
$ws.Range("N73").Value = -9023.25
$ws.Range("H132").Value = 944.16364
$ws.Range("I132").Value = 630.02563
$ws.Range("J132").Value = 1709.875
$ws.Range("K132").Value = 1890.07689
$ws.Range("L132").Value = 5129.625
$ws.Range("M132").Value = 639.9231100000002
$ws.Range("N132").Value = -10189.625
$ws.Range("H137").Value = 1432.7188
$ws.Range("I137").Value = 1290.6296
$ws.Range("J137").Value = 2200
$ws.Range("K137").Value = 3871.8888
$ws.Range("L137").Value = 6600
$ws.Range("M137").Value = -1321.8888
$ws.Range("N137").Value = -11700

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1187.4032
$ws.Range("I61").Value = 971.58
$ws.Range("J61").Value = 2086.6667
$ws.Range("K61").Value = 971.58
$ws.Range("L61").Value = 2086.6667
$ws.Range("M61").Value = -759.58
$ws.Range("N61").Value = -2510.6667
$ws.Range("H74").Value = 1379.1666
$ws.Range("I74").Value = 1238.0435
$ws.Range("J74").Value = 1842.8572
$ws.Range("K74").Value = 1238.0435
$ws.Range("L74").Value = 1842.8572
$ws.Range("M74").Value = -364.0435
$ws.Range("N74").Value = -3590.8572
$ws.Range("H77").Value = 1379.1666
$ws.Range("I77").Value = 1238.0435
$ws.Range("J77").Value = 1842.8572
$ws.Range("K77").Value = 6190.2175
$ws.Range("L77").Value = 9214.286
$ws.Range("M77").Value = -1822.2175
$ws.Range("N77").Value = -17950.286
$ws.Range("H132").Value = 2520.9788
$ws.Range("I132").Value = 1539.3939
$ws.Range("J132").Value = 4834.7144
$ws.Range("K132").Value = 4618.1817
$ws.Range("L132").Value = 14504.1432
$ws.Range("M132").Value = -2088.1817
$ws.Range("N132").Value = -19564.1432
$ws.Range("H136").Value = 1187.4032
$ws.Range("I136").Value = 971.58
$ws.Range("J136").Value = 2086.6667
$ws.Range("K136").Value = 2914.74
$ws.Range("L136").Value = 6260.000100000001
$ws.Range("M136").Value = -364.7400000000002
$ws.Range("N136").Value = -11360.0001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 736.11536
$ws.Range("I22").Value = 590.58826
$ws.Range("J22").Value = 1011
$ws.Range("K22").Value = 590.58826
$ws.Range("L22").Value = 1011
$ws.Range("M22").Value = -417.58826
$ws.Range("N22").Value = -1357
$ws.Range("H134").Value = 1575.9578
$ws.Range("I134").Value = 1351.3396
$ws.Range("J134").Value = 2237.3333
$ws.Range("K134").Value = 4054.0188
$ws.Range("L134").Value = 6711.999899999999
$ws.Range("M134").Value = -1519.0188
$ws.Range("N134").Value = -11781.9999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1258.6129
$ws.Range("I58").Value = 917.9524
$ws.Range("J58").Value = 1974
$ws.Range("K58").Value = 917.9524
$ws.Range("L58").Value = 1974
$ws.Range("M58").Value = -714.9524
$ws.Range("N58").Value = -2380
$ws.Range("H94").Value = 3347.3704
$ws.Range("I94").Value = 4364
$ws.Range("J94").Value = 2839.0557
$ws.Range("K94").Value = 4364
$ws.Range("L94").Value = 2839.0557
$ws.Range("M94").Value = -3913
$ws.Range("N94").Value = -3741.0557
$ws.Range("H136").Value = 1258.6129
$ws.Range("I136").Value = 917.9524
$ws.Range("J136").Value = 1974
$ws.Range("K136").Value = 2753.8572
$ws.Range("L136").Value = 5922
$ws.Range("M136").Value = -203.8571999999999
$ws.Range("N136").Value = -11022

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H32").Value = 3311.111
$ws.Range("I32").Value = 1000
$ws.Range("J32").Value = 3600
$ws.Range("K32").Value = 3000
$ws.Range("L32").Value = 10800
$ws.Range("M32").Value = -2717
$ws.Range("N32").Value = -11366
$ws.Range("H49").Value = 387.5
$ws.Range("I49").Value = 500
$ws.Range("J49").Value = 50
$ws.Range("K49").Value = 1500
$ws.Range("L49").Value = 150
$ws.Range("M49").Value = -1344
$ws.Range("N49").Value = -462
$ws.Range("H58").Value = 2695.3
$ws.Range("I58").Value = 1200
$ws.Range("J58").Value = 2774
$ws.Range("K58").Value = 3600
$ws.Range("L58").Value = 8322
$ws.Range("M58").Value = -3472
$ws.Range("N58").Value = -8578
$ws.Range("H64").Value = 3089.25
$ws.Range("I64").Value = 800
$ws.Range("J64").Value = 3416.2856
$ws.Range("K64").Value = 2400
$ws.Range("L64").Value = 10248.8568
$ws.Range("M64").Value = -2130
$ws.Range("N64").Value = -10788.8568
$ws.Range("H67").Value = 3089.25
$ws.Range("I67").Value = 800
$ws.Range("J67").Value = 3416.2856
$ws.Range("K67").Value = 2400
$ws.Range("L67").Value = 10248.8568
$ws.Range("M67").Value = -1464
$ws.Range("N67").Value = -12120.8568
$ws.Range("H76").Value = 3052
$ws.Range("I76").Value = 2104
$ws.Range("J76").Value = 4000
$ws.Range("K76").Value = 6312
$ws.Range("L76").Value = 12000
$ws.Range("M76").Value = -5929
$ws.Range("N76").Value = -12766
$ws.Range("H79").Value = 3052
$ws.Range("I79").Value = 2104
$ws.Range("J79").Value = 4000
$ws.Range("K79").Value = 6312
$ws.Range("L79").Value = 12000
$ws.Range("M79").Value = -4986
$ws.Range("N79").Value = -14652
$ws.Range("H105").Value = 6851
$ws.Range("I105").Value = 5526
$ws.Range("J105").Value = 6939.3335
$ws.Range("K105").Value = 16578
$ws.Range("L105").Value = 20818.0005
$ws.Range("M105").Value = -13957
$ws.Range("N105").Value = -26060.0005
$ws.Range("H108").Value = 67175.664
$ws.Range("I108").Value = 67175.664
$ws.Range("J108").Value = 0
$ws.Range("K108").Value = 201526.992
$ws.Range("L108").Value = 0
$ws.Range("M108").Value = -198646.992
$ws.Range("H124").Value = 5874.5454
$ws.Range("I124").Value = 2860
$ws.Range("J124").Value = 6544.4443
$ws.Range("K124").Value = 8580
$ws.Range("L124").Value = 19633.3329
$ws.Range("M124").Value = -3670
$ws.Range("N124").Value = -29453.3329
$ws.Range("H136").Value = 12436.333
$ws.Range("I136").Value = 14699.143
$ws.Range("J136").Value = 4516.5
$ws.Range("K136").Value = 44097.429
$ws.Range("L136").Value = 13549.5
$ws.Range("M136").Value = -38997.429
$ws.Range("N136").Value = -23749.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2803.4092
$ws.Range("I80").Value = 2453.5715
$ws.Range("J80").Value = 2966.6667
$ws.Range("K80").Value = 2453.5715
$ws.Range("L80").Value = 2966.6667
$ws.Range("M80").Value = -1455.5715
$ws.Range("N80").Value = -4962.6667
$ws.Range("H83").Value = 2803.4092
$ws.Range("I83").Value = 2453.5715
$ws.Range("J83").Value = 2966.6667
$ws.Range("K83").Value = 12267.8575
$ws.Range("L83").Value = 14833.3335
$ws.Range("M83").Value = -7275.8575
$ws.Range("N83").Value = -24817.3335
$ws.Range("H93").Value = 9251
$ws.Range("I93").Value = 0
$ws.Range("J93").Value = 9251
$ws.Range("K93").Value = 0
$ws.Range("L93").Value = 9251
$ws.Range("N93").Value = -12995
$ws.Range("H122").Value = 1622652.1
$ws.Range("I122").Value = 2494756.8
$ws.Range("J122").Value = 3029.1428
$ws.Range("K122").Value = 7484270.399999999
$ws.Range("L122").Value = 9087.428400000001
$ws.Range("M122").Value = -7481820.399999999
$ws.Range("N122").Value = -13987.4284
$ws.Range("H123").Value = 22097.625
$ws.Range("I123").Value = 0
$ws.Range("J123").Value = 22097.625
$ws.Range("K123").Value = 0
$ws.Range("L123").Value = 22097.625
$ws.Range("N123").Value = -26997.625

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 4275217
$ws.Range("I22").Value = 13889507
$ws.Range("J22").Value = 2199.6667
$ws.Range("K22").Value = 13889507
$ws.Range("L22").Value = 2199.6667
$ws.Range("M22").Value = -13889212
$ws.Range("N22").Value = -2789.6667
$ws.Range("H27").Value = 4275217
$ws.Range("I27").Value = 13889507
$ws.Range("J27").Value = 2199.6667
$ws.Range("K27").Value = 13889507
$ws.Range("L27").Value = 2199.6667
$ws.Range("M27").Value = -13889400
$ws.Range("N27").Value = -2413.6667
$ws.Range("H46").Value = 19608794
$ws.Range("I46").Value = 41667388
$ws.Range("J46").Value = 1155.5555
$ws.Range("K46").Value = 41667388
$ws.Range("L46").Value = 1155.5555
$ws.Range("M46").Value = -41667200
$ws.Range("N46").Value = -1531.5555
$ws.Range("H55").Value = 78947610
$ws.Range("I55").Value = 90909300
$ws.Range("J55").Value = 62500270
$ws.Range("K55").Value = 90909300
$ws.Range("L55").Value = 62500270
$ws.Range("M55").Value = -90909127
$ws.Range("N55").Value = -62500616

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H47").Value = 0
$ws.Range("I47").Value = 0
$ws.Range("J47").Value = 0
$ws.Range("K47").Value = 0
$ws.Range("L47").Value = 0
$ws.Range("N47").ClearContents()
$ws.Range("H132").Value = 1141.6863
$ws.Range("I132").Value = 835.575
$ws.Range("J132").Value = 2254.818
$ws.Range("K132").Value = 2506.725
$ws.Range("L132").Value = 6764.454000000001
$ws.Range("M132").Value = 23.27499999999964
$ws.Range("N132").Value = -11824.454
